$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values in column D (rows 2-33), minor revisions from source data update
$dUpdates = @{
    2  = 33924
    3  = 34448
    4  = 34742
    5  = 34815
    6  = 34867
    7  = 35031
    8  = 35002
    9  = 35406
    10 = 35500
    11 = 35965
    12 = 35830
    13 = 36180
    14 = 36418
    15 = 36341
    16 = 36473
    17 = 36569
    18 = 36268
    19 = 36646
    20 = 37303
    21 = 37635
    22 = 38055
    23 = 38386
    24 = 38272
    25 = 38706
    26 = 38599
    27 = 39019
    28 = 39331
    29 = 37818
    30 = 38552
    31 = 33665
    32 = 35493
    33 = 37788
}

foreach ($row in $dUpdates.Keys) {
    $ws.Cells.Item($row, 4).Value = $dUpdates[$row]
}

# Row 34 gets revised values for B, C, D
$ws.Cells.Item(34, 2).Value = 56433
$ws.Cells.Item(34, 3).Value = 38076
$ws.Cells.Item(34, 4).Value = 39060

# New row 35 for the latest quarter (01-04-2021).
# Force text formatting first so Excel stores the literal string (shared
# string) instead of auto-converting the date-looking text to a serial
# date number, then restore the default (unstyled) cell format.
$ws.Cells.Item(35, 1).NumberFormat = "@"
$ws.Cells.Item(35, 1).Value = "01-04-2021"
$ws.Cells.Item(35, 1).Style = "Normal"
$ws.Cells.Item(35, 2).Value = 57427
$ws.Cells.Item(35, 3).Value = 39677
$ws.Cells.Item(35, 4).Value = 39460
